# Update "想去人数" (number of people interested) values on the
# "展览" and "全部类型" worksheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 764
$wsExhibit.Range("F4").Value = 265
$wsExhibit.Range("F5").Value = 793
$wsExhibit.Range("F6").Value = 1926
$wsExhibit.Range("F7").Value = 171

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 764
$wsAll.Range("F4").Value = 265
$wsAll.Range("F7").Value = 793
$wsAll.Range("F8").Value = 1926
$wsAll.Range("F10").Value = 171
